$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I33").Value = 195.4
$ws.Range("K33").Value = 195.4
$ws.Range("M33").Value = 33.59999999999999
$ws.Range("H53").Value = 2276.875
$ws.Range("I53").Value = 3531.2666
$ws.Range("J53").Value = 186.22223
$ws.Range("K53").Value = 3531.2666
$ws.Range("L53").Value = 186.22223
$ws.Range("M53").Value = -2894.2666
$ws.Range("N53").Value = -1460.22223
$ws.Range("H86").Value = 1301658.8
$ws.Range("I86").Value = 3297.2727
$ws.Range("K86").Value = 3297.2727
$ws.Range("M86").Value = -2174.2727
$ws.Range("H89").Value = 1301658.8
$ws.Range("I89").Value = 3297.2727
$ws.Range("K89").Value = 16486.3635
$ws.Range("M89").Value = -10870.3635
$ws.Range("H98").Value = 1540.5227
$ws.Range("J98").Value = 1532.6
$ws.Range("L98").Value = 1532.6
$ws.Range("N98").Value = -4528.6
$ws.Range("H103").Value = 441.5
$ws.Range("J103").Value = 416.33334
$ws.Range("L103").Value = 1249.00002
$ws.Range("N103").Value = -2421.00002
$ws.Range("H112").Value = 1524.65
$ws.Range("J112").Value = 1553.7632
$ws.Range("L112").Value = 4661.2896
$ws.Range("N112").Value = -6877.2896
$ws.Range("H122").Value = 1540.5227
$ws.Range("J122").Value = 1532.6
$ws.Range("L122").Value = 4597.799999999999
$ws.Range("N122").Value = -9497.799999999999
$ws.Range("H137").Value = 2253.077
$ws.Range("I137").Value = 2693.5217
$ws.Range("J137").Value = 1619.9375
$ws.Range("K137").Value = 8080.5651
$ws.Range("L137").Value = 4859.8125
$ws.Range("M137").Value = -5530.5651
$ws.Range("N137").Value = -9959.8125
$ws.Range("H138").Value = 13516318
$ws.Range("J138").Value = 3106.39
$ws.Range("L138").Value = 9319.17
$ws.Range("N138").Value = -19599.17

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1950.5714
$ws.Range("I45").Value = 2025.8334
$ws.Range("K45").Value = 2025.8334
$ws.Range("M45").Value = -1648.8334
$ws.Range("H97").Value = 699.5
$ws.Range("I97").Value = 699.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 699.5
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 5878.5
$ws.Range("I122").Value = 6062.5
$ws.Range("J122").Value = 4406.5
$ws.Range("K122").Value = 18187.5
$ws.Range("L122").Value = 13219.5
$ws.Range("M122").Value = -15737.5
$ws.Range("N122").Value = -18119.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1999.8334
$ws.Range("I105").Value = 1999.8334
$ws.Range("K105").Value = 1999.8334
$ws.Range("M105").Value = -252.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1411.091
$ws.Range("I132").Value = 1252.2
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3756.6
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1226.6
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 1906.4348
$ws.Range("I134").Value = 1876.4286
$ws.Range("K134").Value = 5629.2858
$ws.Range("M134").Value = -3094.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1359.75
$ws.Range("I14").Value = 1359.75
$ws.Range("K14").Value = 4079.25
$ws.Range("M14").Value = -3906.25
$ws.Range("H121").Value = 2594.9375
$ws.Range("J121").Value = 2917.7856
$ws.Range("L121").Value = 8753.356800000001
$ws.Range("N121").Value = -11373.3568
$ws.Range("H122").Value = 3087.9285
$ws.Range("I122").Value = 3916
$ws.Range("K122").Value = 35244
$ws.Range("M122").Value = -32794
$ws.Range("H132").Value = 2199.5
$ws.Range("J132").Value = 2199
$ws.Range("L132").Value = 19791
$ws.Range("N132").Value = -24851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5700.5
$ws.Range("I126").Value = 4286.4287
$ws.Range("K126").Value = 12859.2861
$ws.Range("M126").Value = -10389.2861
$ws.Range("H136").Value = 56817.473
$ws.Range("J136").Value = 56817.473
$ws.Range("L136").Value = 170452.419
$ws.Range("N136").Value = -175552.419

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 626.3182
$ws.Range("J16").Value = 244.5
$ws.Range("L16").Value = 244.5
$ws.Range("N16").Value = -584.5
$ws.Range("H22").Value = 3795.3845
$ws.Range("J22").Value = 6253.143
$ws.Range("L22").Value = 6253.143
$ws.Range("N22").Value = -6843.143
$ws.Range("H27").Value = 3795.3845
$ws.Range("J27").Value = 6253.143
$ws.Range("L27").Value = 6253.143
$ws.Range("N27").Value = -6467.143
$ws.Range("H46").Value = 7576617.5
$ws.Range("I46").Value = 8334185
$ws.Range("K46").Value = 8334185
$ws.Range("M46").Value = -8333997
$ws.Range("H68").Value = 4857.6665
$ws.Range("J68").Value = 10328.75
$ws.Range("L68").Value = 10328.75
$ws.Range("N68").Value = -11826.75
$ws.Range("H71").Value = 4857.6665
$ws.Range("J71").Value = 10328.75
$ws.Range("L71").Value = 51643.75
$ws.Range("N71").Value = -59131.75
$ws.Range("H93").Value = 2262.9678
$ws.Range("I93").Value = 1396.5217
$ws.Range("K93").Value = 1396.5217
$ws.Range("M93").Value = -148.5217
$ws.Range("H132").Value = 2851.9893
$ws.Range("I132").Value = 2859.7654
$ws.Range("J132").Value = 2799.5
$ws.Range("K132").Value = 8579.296200000001
$ws.Range("L132").Value = 8398.5
$ws.Range("M132").Value = -6049.296200000001
$ws.Range("N132").Value = -13458.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 21103.4
$ws.Range("I41").Value = 21681
$ws.Range("J41").Value = 20718.334
$ws.Range("K41").Value = 21681
$ws.Range("L41").Value = 20718.334
$ws.Range("M41").Value = -21291
$ws.Range("N41").Value = -21498.334
$ws.Range("H62").Value = 35153.918
$ws.Range("I62").Value = 22184.9
$ws.Range("J62").Value = 99999
$ws.Range("K62").Value = 22184.9
$ws.Range("L62").Value = 99999
$ws.Range("M62").Value = -21560.9
$ws.Range("N62").Value = -101247
$ws.Range("H65").Value = 35153.918
$ws.Range("I65").Value = 22184.9
$ws.Range("J65").Value = 99999
$ws.Range("K65").Value = 110924.5
$ws.Range("L65").Value = 499995
$ws.Range("M65").Value = -107804.5
$ws.Range("N65").Value = -506235
$ws.Range("H132").Value = 823.7368
$ws.Range("I132").Value = 861.8823
$ws.Range("K132").Value = 2585.6469
$ws.Range("M132").Value = -55.64689999999973
$ws.Range("H136").Value = 3631.311
$ws.Range("I136").Value = 3823.9143
$ws.Range("K136").Value = 11471.7429
$ws.Range("M136").Value = -8921.742899999999
